# Update "想去人数" (want-to-go count) figures in the F column across the
# relevant worksheets, per the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibition.Range("F2").Value = 6700
$wsExhibition.Range("F4").Value = 423
$wsExhibition.Range("F12").Value = 171
$wsExhibition.Range("F13").Value = 402
$wsExhibition.Range("F14").Value = 1
$wsExhibition.Range("F15").Value = 1450
$wsExhibition.Range("F16").Value = 18
$wsExhibition.Range("F17").Value = 3345
$wsExhibition.Range("F19").Value = 222
$wsExhibition.Range("F20").Value = 2
$wsExhibition.Range("F21").Value = 1995
$wsExhibition.Range("F22").Value = 107
$wsExhibition.Range("F23").Value = 28
$wsExhibition.Range("F25").Value = 130

# 演出 (sheet2)
$wsPerformance.Range("F2").Value = 4

# 全部类型 (sheet4)
$wsAllTypes.Range("F2").Value = 6700
$wsAllTypes.Range("F4").Value = 423
$wsAllTypes.Range("F7").Value = 4
$wsAllTypes.Range("F13").Value = 171
$wsAllTypes.Range("F14").Value = 402
$wsAllTypes.Range("F15").Value = 1
$wsAllTypes.Range("F16").Value = 1450
$wsAllTypes.Range("F17").Value = 18
$wsAllTypes.Range("F18").Value = 3345
$wsAllTypes.Range("F20").Value = 222
$wsAllTypes.Range("F21").Value = 2
$wsAllTypes.Range("F22").Value = 1995
$wsAllTypes.Range("F23").Value = 107
$wsAllTypes.Range("F24").Value = 28
$wsAllTypes.Range("F26").Value = 130
